# Commit: "removed unwanted projects, added new code for Excel."
# - Insert a new blank worksheet ("Sheet1") in front of the existing
#   "Info_Table" sheet, and make it the active/selected sheet.
# - On "Info_Table", write a bold header "Your header" into cell G1 and
#   leave that cell selected (via a shared string + bold font style).

$wb = $excel.ActiveWorkbook

# Adding a worksheet inserts it directly before the currently active
# sheet and makes the new sheet active - this gives us "Sheet1" as the
# first, selected tab and "Info_Table" pushed to second position.
# (Look up "Info_Table" again afterwards - sheet references in this
# runtime resolve by position, and the insert shifts its index.)
$newSheet = $wb.Worksheets.Add()
$infoSheet = $wb.Worksheets.Item("Info_Table")

# Add the bold header value on the Info_Table sheet.
$headerCell = $infoSheet.Range("G1")
$headerCell.Value = "Your header"
$headerCell.Font.Bold = $true
$headerCell.Select() | Out-Null

# Restore the new "Sheet1" as the active/selected sheet.
$newSheet.Select() | Out-Null
